# DataFile_LiveSLR_Testing.xlsx - add an "ExpectedSourceTemplateFile" column
# A new column is inserted before the existing column I ("Report-" filename
# column, which shifts right to become column J), and populated with a
# header + one data value describing the expected source template file used
# for the "content comparison btw Excel reports" test.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at I; everything from old column I onward
# (I..M) shifts one column to the right (new J..N).
$ws.Columns("I:I").Insert()

# The newly inserted column inherits the width of its left neighbour (H)
# in the original file, same as Excel's own "Insert Column" behaviour.
$ws.Columns("I:I").ColumnWidth = $ws.Columns("H:H").ColumnWidth

# Populate the new column's header (row 1) and its single data value (row 2).
$ws.Cells.Item(1, 9).Value = "ExpectedSourceTemplateFile"
$ws.Cells.Item(2, 9).Value = "\Testdata\Templates\SLRReport_SourceData\Expected_Source_Data_Manipulated.xlsx"

# Reflect the author's final selection/view state (active cell on the new
# value cell I2).
[void]$ws.Range("I2").Select()
